# Auto-generated Excel COM-interop script to apply market data updates
# to the Kraken_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 500
$ws.Range("K2").Value = 500
$ws.Range("M2").Value = -387
$ws.Range("H9").Value = 408.16666
$ws.Range("I9").Value = 408.16666
$ws.Range("K9").Value = 408.16666
$ws.Range("M9").Value = -239.16666
$ws.Range("H19").Value = 621.1667
$ws.Range("I19").Value = 606.75
$ws.Range("J19").Value = 650
$ws.Range("K19").Value = 606.75
$ws.Range("L19").Value = 650
$ws.Range("M19").Value = -431.75
$ws.Range("N19").Value = -1000
$ws.Range("H88").Value = 9136
$ws.Range("I88").Value = 9999
$ws.Range("K88").Value = 9999
$ws.Range("M88").Value = -9593
$ws.Range("H91").Value = 9136
$ws.Range("I91").Value = 9999
$ws.Range("K91").Value = 9999
$ws.Range("M91").Value = -8595
$ws.Range("H107").Value = 1015.2143
$ws.Range("I107").Value = 771.3
$ws.Range("J107").Value = 1625
$ws.Range("K107").Value = 771.3
$ws.Range("L107").Value = 1625
$ws.Range("M107").Value = 1148.7
$ws.Range("N107").Value = -5465
$ws.Range("H111").Value = 471.875
$ws.Range("I111").Value = 345.5
$ws.Range("J111").Value = 598.25
$ws.Range("K111").Value = 1036.5
$ws.Range("L111").Value = 1794.75
$ws.Range("M111").Value = 2030.5
$ws.Range("N111").Value = -7928.75
$ws.Range("H135").Value = 2431.8
$ws.Range("I135").Value = 2489.75
$ws.Range("J135").Value = 2200
$ws.Range("K135").Value = 22407.75
$ws.Range("L135").Value = 19800
$ws.Range("M135").Value = -19872.75
$ws.Range("N135").Value = -24870
$ws.Range("H138").Value = 3568.1
$ws.Range("I138").Value = 1126.3334
$ws.Range("K138").Value = 3379.0002
$ws.Range("M138").Value = 1760.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4875
$ws.Range("I61").Value = 4875
$ws.Range("K61").Value = 4875
$ws.Range("M61").Value = -4663
$ws.Range("H74").Value = 3518.1428
$ws.Range("I74").Value = 3938
$ws.Range("J74").Value = 999
$ws.Range("K74").Value = 3938
$ws.Range("L74").Value = 999
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -2747
$ws.Range("H77").Value = 3518.1428
$ws.Range("I77").Value = 3938
$ws.Range("J77").Value = 999
$ws.Range("K77").Value = 19690
$ws.Range("L77").Value = 4995
$ws.Range("M77").Value = -15322
$ws.Range("N77").Value = -13731
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("N122").Value = 0
$ws.Range("H136").Value = 4875
$ws.Range("I136").Value = 4875
$ws.Range("K136").Value = 14625
$ws.Range("M136").Value = -12075
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 9084.6
$ws.Range("I107").Value = 4365.5454
$ws.Range("K107").Value = 4365.5454
$ws.Range("M107").Value = -2445.5454
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6733.0557
$ws.Range("I31").Value = 5800.2144
$ws.Range("K31").Value = 5800.2144
$ws.Range("M31").Value = -5505.2144
$ws.Range("H34").Value = 6733.0557
$ws.Range("I34").Value = 5800.2144
$ws.Range("K34").Value = 5800.2144
$ws.Range("M34").Value = -5598.2144
$ws.Range("H86").Value = 7999.6665
$ws.Range("I86").Value = 3999.5
$ws.Range("K86").Value = 3999.5
$ws.Range("M86").Value = -2876.5
$ws.Range("H89").Value = 7999.6665
$ws.Range("I89").Value = 3999.5
$ws.Range("K89").Value = 19997.5
$ws.Range("M89").Value = -14381.5
$ws.Range("H99").Value = 1166.6666
$ws.Range("J99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496
$ws.Range("H107").Value = 669.6667
$ws.Range("I107").Value = 605.5
$ws.Range("K107").Value = 605.5
$ws.Range("M107").Value = 1314.5
$ws.Range("H122").Value = 1293.25
$ws.Range("I122").Value = 1293.25
$ws.Range("K122").Value = 3879.75
$ws.Range("M122").Value = -1429.75
$ws.Range("H126").Value = 1166.6666
$ws.Range("J126").Value = 1500
$ws.Range("L126").Value = 4500
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3083.25
$ws.Range("I5").Value = 3533.3333
$ws.Range("J5").Value = 2933.2222
$ws.Range("K5").Value = 10599.9999
$ws.Range("L5").Value = 8799.6666
$ws.Range("M5").Value = -10487.9999
$ws.Range("N5").Value = -9023.6666
$ws.Range("H12").Value = 45.8
$ws.Range("I12").Value = 72.333336
$ws.Range("J12").Value = 34.42857
$ws.Range("K12").Value = 217.000008
$ws.Range("L12").Value = 103.28571
$ws.Range("M12").Value = -44.00000800000001
$ws.Range("N12").Value = -449.28571
$ws.Range("H68").Value = 1874.5
$ws.Range("J68").Value = 1832.6666
$ws.Range("L68").Value = 5497.9998
$ws.Range("N68").Value = -7119.9998
$ws.Range("H71").Value = 1874.5
$ws.Range("J71").Value = 1832.6666
$ws.Range("L71").Value = 16493.9994
$ws.Range("N71").Value = -24605.9994
$ws.Range("H121").Value = 1258.9166
$ws.Range("I121").Value = 575
$ws.Range("J121").Value = 1486.8889
$ws.Range("K121").Value = 1725
$ws.Range("L121").Value = 4460.6667
$ws.Range("M121").Value = -415
$ws.Range("N121").Value = -7080.6667
$ws.Range("H135").Value = 3083.25
$ws.Range("I135").Value = 3533.3333
$ws.Range("J135").Value = 2933.2222
$ws.Range("K135").Value = 31799.9997
$ws.Range("L135").Value = 26398.9998
$ws.Range("M135").Value = -29264.9997
$ws.Range("N135").Value = -31468.9998
$ws.Range("H137").Value = 1267.375
$ws.Range("I137").Value = 1055.2
$ws.Range("K137").Value = 3165.6
$ws.Range("M137").Value = 1934.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3410.2727
$ws.Range("I102").Value = 3745.7144
$ws.Range("K102").Value = 3745.7144
$ws.Range("M102").Value = -2123.7144
$ws.Range("H113").Value = 1070
$ws.Range("I113").Value = 940
$ws.Range("K113").Value = 940
$ws.Range("M113").Value = 1230
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 6074.8335
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 1100
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 1100
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -1476
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("L67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3999.8572
$ws.Range("I126").Value = 3666.5
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 10999.5
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -8529.5
$ws.Range("N126").Value = -22940
